# FR-GTH-108 DESCUENTO EXAMENES MEDICOS — date line fix.
#
# The greeting/date line currently reads:
#   Medellín, {{ fecha_dia }} de  {{ fecha_mes }} de {{ fecha_año }}
#
# It must read:
#   Medellín, {{ fecha_dia }} del mes número  {{ fecha_mes }} de {{ fecha_año }}
#
# i.e. the run that only contains "de  " (note: two trailing spaces) turns
# into "del mes número  " (the two trailing spaces before the next "{{" are
# preserved). That exact "de  " text is unique in the whole document, so a
# single, case-sensitive Find/Replace unambiguously targets it without
# touching the "de " (single trailing space) run used later in the same
# sentence.

$d = $word.ActiveDocument

$searchText  = "de  "
$replaceText = "del mes número  "

$found = $d.Content.Find.Execute(
    $searchText,   # FindText
    $true,         # MatchCase
    $false,        # MatchWholeWord
    $false,        # MatchWildcards
    $false,        # MatchSoundsLike
    $false,        # MatchAllWordForms
    $true,         # Forward
    1,             # Wrap (wdFindContinue)
    $false,        # Format
    $replaceText,  # ReplaceWith
    2              # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Expected text 'de  ' was not found in the document."
}

Write-Output "Replace executed: $found"
Write-Output "New date line: $($d.Paragraphs(4).Range.Text)"
